$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab from "op2" to "wong3"
$ws.Name = "wong3"

# Row 2 - Weight.ElevatorSystem.Elevator.ASRS_2
$ws.Range("E2").Value = 18
$ws.Range("F2").Value = 18
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 393

# Row 4 - Empty.ElevatorSystem.Environment.AOIU_1
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 393

# Row 8 - Empty.ElevatorSystem.Environment.AOIU_5
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 3
$ws.Range("E8").Value = 21
$ws.Range("F8").Value = 21

# Row 12 - Empty.PL_Interface_impl.AOIU_1
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 24

# Row 13 - Weight.ElevatorSystem.Elevator.ODL_2
$ws.Range("E13").Value = 9
$ws.Range("J13").Value = 12
$ws.Range("K13").Value = 409
